$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Files")

# Strict / not-strict scan id collation: the DICOM series-instance suffix
# ("_3", "_6168", ...) is dropped from the computed scan id in column I
# (DICOM:SeriesNumber), so duplicate scans collapse onto the same id.
$ws.Cells.Item(2, 9).Value = "002304_CT1"
$ws.Cells.Item(3, 9).Value = "002304_CT1"
$ws.Cells.Item(4, 9).Value = "002304_CT1"
$ws.Cells.Item(5, 9).Value = "002304_CT1"
$ws.Cells.Item(6, 9).Value = "002304_CT1"
$ws.Cells.Item(7, 9).Value = "397829_CT1"
$ws.Cells.Item(8, 9).Value = "397829_CT2"
$ws.Cells.Item(9, 9).Value = "397829_CT3"
$ws.Cells.Item(10, 9).Value = "038945_CT1"

# Move the active selection as it was left in the saved workbook.
$ws.Range("M2").Select()
